{"js": "// The edit:\n//   1) Remove the empty (blank) heading-style paragraph that sits between the\n//      document title (\"Part 2a: Testing Plan & Test Cases\") and the\n//      \"Testing Strategy\" heading.\n//   2) Merge the trailing \".\" run back into the preceding run of the\n//      \"Finally, a combined test stub ...\" paragraph so the sentence ends up\n//      as a single run (\".\", previously its own run, is joined to the text\n//      before it).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Edit 1: delete the blank paragraph right after the title paragraph ---\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"Part 2a: Testing Plan & Test Cases\") {\n    const next = items[i + 1];\n    if (next) {\n      next.load(\"text\");\n      await context.sync();\n      if (next.text.trim() === \"\") {\n        next.delete();\n        await context.sync();\n      }\n    }\n    break;\n  }\n}\n\n// --- Edit 2: merge the split runs of the \"Finally, a combined...\" paragraph ---\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n  const para = paragraphs2.items[i];\n  if (para.text.indexOf(\"Finally, a combined test stub\") === 0) {\n    const fullText = para.text;\n    // Re-writing the paragraph's text collapses it back down to a single run,\n    // joining the separate \".\" run into the sentence that precedes it.\n    para.insertText(fullText, \"Replace\");\n    await context.sync();\n    break;\n  }\n}\n", "ps1": "# The edit:\n#   1) Remove the empty (blank) heading-style paragraph that sits between the\n#      document title (\"Part 2a: Testing Plan & Test Cases\") and the\n#      \"Testing Strategy\" heading.\n#   2) Merge the trailing \".\" run back into the preceding run of the\n#      \"Finally, a combined test stub ...\" paragraph so the sentence ends up\n#      as a single run.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: delete the blank paragraph right after the title paragraph ---\n$titleIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $ptext = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($ptext -eq \"Part 2a: Testing Plan & Test Cases\") {\n        $titleIdx = $i\n        break\n    }\n}\n\nif ($titleIdx -gt 0) {\n    $nextPara = $d.Paragraphs.Item($titleIdx + 1)\n    if ($nextPara.Range.Text.Trim() -eq \"\") {\n        $nextPara.Range.Delete()\n    }\n}\n\n# --- Edit 2: merge the split runs of the \"Finally, a combined...\" paragraph ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"movement).\"\n$find.Replacement.Text = \"movement).\"\n$find.Execute([ref]\"movement).\", $false, $false, $false, $false, $false, $true, 1, $false, \"movement).\", 2)\n"}
